$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Status text update: "Ready for handoff" -> "Handed back: in sync with en-us"
# (appears in column B, rows 2 & 3, on both the zh-cn and de-de sheets)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-us"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("B3").Value = $newStatus
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("B3").Value = $newStatus

# ---------------------------------------------------------------------------
# zh-cn sheet: populate "Latest Target File" (E) / "Latest Handback File" (F)
# columns for the two tracked files, and record the handback timestamp in
# "Latest Handback DateTime" (G).
# ---------------------------------------------------------------------------
$wsZh.Range("E2").Value = "13e04a81-d56b-4947-8f13-c5477b440c36.md"
$wsZh.Range("F2").Value = "13e04a81-d56b-4947-8f13-c5477b440c36.64c709783fa5e5598193b310a8574c4499ab0bf1.zh-cn.xlf"
$wsZh.Range("G2").Value = "2016-01-08 18:02:41"

$wsZh.Range("E3").Value = "1e023323-5838-45bb-bfa0-693e87d12526.md"
$wsZh.Range("F3").Value = "1e023323-5838-45bb-bfa0-693e87d12526.839f7e752aa1ad3605eaf4ec4e62d3d74e11bdfc.zh-cn.xlf"
$wsZh.Range("G3").Value = "2016-01-08 18:02:41"

# ---------------------------------------------------------------------------
# de-de sheet: same treatment
# ---------------------------------------------------------------------------
$wsDe.Range("E2").Value = "13e04a81-d56b-4947-8f13-c5477b440c36.md"
$wsDe.Range("F2").Value = "13e04a81-d56b-4947-8f13-c5477b440c36.64c709783fa5e5598193b310a8574c4499ab0bf1.de-de.xlf"
$wsDe.Range("G2").Value = "2016-01-08 18:02:57"

$wsDe.Range("E3").Value = "1e023323-5838-45bb-bfa0-693e87d12526.md"
$wsDe.Range("F3").Value = "1e023323-5838-45bb-bfa0-693e87d12526.839f7e752aa1ad3605eaf4ec4e62d3d74e11bdfc.de-de.xlf"
$wsDe.Range("G3").Value = "2016-01-08 18:02:57"

# ---------------------------------------------------------------------------
# Rebuild the hyperlinks on both sheets so the new E/F columns get linked the
# same way their A/C counterparts are, and so the relationship ids come out
# in a stable, predictable order: A2, C2, E2, F2, A3, C3, E3, F3, A4.
# ---------------------------------------------------------------------------
$mdTarget1 = "https://github.com/OpenLocalizationTest/oltest/blob/7f9255525c7fd8567062d89745b78e3169a40d9b/e2e/13e04a81-d56b-4947-8f13-c5477b440c36.md"
$mdTarget2 = "https://github.com/OpenLocalizationTest/oltest/blob/7f9255525c7fd8567062d89745b78e3169a40d9b/e2e/1e023323-5838-45bb-bfa0-693e87d12526.md"
$configTarget = "https://github.com/OpenLocalizationTest/oltest/blob/7f9255525c7fd8567062d89745b78e3169a40d9b/.localization-config"

$zhXlf2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5cb668b97a23106d5ae428771d0bfe4e8e781817/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/13e04a81-d56b-4947-8f13-c5477b440c36.64c709783fa5e5598193b310a8574c4499ab0bf1.zh-cn.xlf"
$zhXlf3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5cb668b97a23106d5ae428771d0bfe4e8e781817/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/1e023323-5838-45bb-bfa0-693e87d12526.839f7e752aa1ad3605eaf4ec4e62d3d74e11bdfc.zh-cn.xlf"

$deXlf2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/339471af15779408e3346b43ed0c36ae88beee5d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/13e04a81-d56b-4947-8f13-c5477b440c36.64c709783fa5e5598193b310a8574c4499ab0bf1.de-de.xlf"
$deXlf3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/339471af15779408e3346b43ed0c36ae88beee5d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/1e023323-5838-45bb-bfa0-693e87d12526.839f7e752aa1ad3605eaf4ec4e62d3d74e11bdfc.de-de.xlf"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdTarget1, "", "", "13e04a81-d56b-4947-8f13-c5477b440c36.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhXlf2, "", "", "13e04a81-d56b-4947-8f13-c5477b440c36.64c709783fa5e5598193b310a8574c4499ab0bf1.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $mdTarget1, "", "", "13e04a81-d56b-4947-8f13-c5477b440c36.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhXlf2, "", "", "13e04a81-d56b-4947-8f13-c5477b440c36.64c709783fa5e5598193b310a8574c4499ab0bf1.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdTarget2, "", "", "1e023323-5838-45bb-bfa0-693e87d12526.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), $zhXlf3, "", "", "1e023323-5838-45bb-bfa0-693e87d12526.839f7e752aa1ad3605eaf4ec4e62d3d74e11bdfc.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), $mdTarget2, "", "", "1e023323-5838-45bb-bfa0-693e87d12526.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhXlf3, "", "", "1e023323-5838-45bb-bfa0-693e87d12526.839f7e752aa1ad3605eaf4ec4e62d3d74e11bdfc.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $configTarget, "", "", ".localization-config")

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdTarget1, "", "", "13e04a81-d56b-4947-8f13-c5477b440c36.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deXlf2, "", "", "13e04a81-d56b-4947-8f13-c5477b440c36.64c709783fa5e5598193b310a8574c4499ab0bf1.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $mdTarget1, "", "", "13e04a81-d56b-4947-8f13-c5477b440c36.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deXlf2, "", "", "13e04a81-d56b-4947-8f13-c5477b440c36.64c709783fa5e5598193b310a8574c4499ab0bf1.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdTarget2, "", "", "1e023323-5838-45bb-bfa0-693e87d12526.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), $deXlf3, "", "", "1e023323-5838-45bb-bfa0-693e87d12526.839f7e752aa1ad3605eaf4ec4e62d3d74e11bdfc.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), $mdTarget2, "", "", "1e023323-5838-45bb-bfa0-693e87d12526.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deXlf3, "", "", "1e023323-5838-45bb-bfa0-693e87d12526.839f7e752aa1ad3605eaf4ec4e62d3d74e11bdfc.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $configTarget, "", "", ".localization-config")
